$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price (D) and volume/change (E) columns with latest scraped values.
# Values that look like plain numbers (single decimal point) are prefixed with a leading
# apostrophe so Excel stores them as text (matching the original inline-string cell type)
# instead of silently coercing them into floating point numbers.

$ws.Range("D2").Value = "66.866.78"
$ws.Range("E2").Value = "  -0.58%  "
$ws.Range("D3").Value = "3.094.91"
$ws.Range("E3").Value = "  -0.66%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'576.54"
$ws.Range("E5").Value = "  -0.60%  "
$ws.Range("D6").Value = "'178.27"
$ws.Range("E6").Value = "  +2.68%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "3.093.07"
$ws.Range("E8").Value = "  -0.52%  "
$ws.Range("D9").Value = "'0.514"
$ws.Range("E9").Value = "  -1.35%  "
$ws.Range("D10").Value = "'6.39"
$ws.Range("E10").Value = "  -2.13%  "
$ws.Range("E11").Value = "  -2.00%  "
$ws.Range("D12").Value = "'0.467"
$ws.Range("E12").Value = "  -2.75%  "
$ws.Range("D13").Value = "'0.0000241"
$ws.Range("E13").Value = "  -3.14%  "
$ws.Range("D14").Value = "'36.09"
$ws.Range("E14").Value = "  -2.22%  "
$ws.Range("E15").Value = "  -0.55%  "
$ws.Range("D16").Value = "3.609.80"
$ws.Range("E16").Value = "  -0.64%  "
$ws.Range("D17").Value = "66.788.06"
$ws.Range("E17").Value = "  -0.65%  "
$ws.Range("D18").Value = "'7.00"
$ws.Range("E18").Value = "  -1.46%  "
$ws.Range("D19").Value = "'16.77"
$ws.Range("E19").Value = "  +0.85%  "
$ws.Range("D20").Value = "3.090.92"
$ws.Range("E20").Value = "  -0.74%  "
$ws.Range("D21").Value = "'481.14"
$ws.Range("E21").Value = "  -2.35%  "
$ws.Range("D22").Value = "'7.75"
$ws.Range("E22").Value = "  -1.96%  "
$ws.Range("D23").Value = "'0.691"
$ws.Range("E23").Value = "  -2.18%  "
$ws.Range("D24").Value = "'83.54"
$ws.Range("E24").Value = "  -0.53%  "
$ws.Range("D25").Value = "'12.67"
$ws.Range("E25").Value = "  -4.14%  "
$ws.Range("E26").Value = "  -2.68%  "
$ws.Range("E27").Value = "  -4.48%  "
$ws.Range("D29").Value = "'7.98"
$ws.Range("E29").Value = "  +0.06%  "
$ws.Range("E30").Value = "  -4.12%  "
$ws.Range("E31").Value = "  -2.79%  "
$ws.Range("D32").Value = "'27.98"
$ws.Range("E32").Value = "  -1.58%  "
$ws.Range("E34").Value = "  -0.15%  "
$ws.Range("D35").Value = "'0.999"
$ws.Range("E35").Value = "  -0.01%  "
$ws.Range("D36").Value = "'48.49"
$ws.Range("E36").Value = "  +2.56%  "
$ws.Range("D37").Value = "'5.60"
$ws.Range("E37").Value = "  -4.78%  "
$ws.Range("D38").Value = "'0.942"
$ws.Range("E38").Value = "  -3.48%  "
$ws.Range("D39").Value = "'0.312"
$ws.Range("E39").Value = "  +0.59%  "
$ws.Range("D40").Value = "'49.00"
$ws.Range("E40").Value = "  -2.11%  "
$ws.Range("E41").Value = "  -2.30%  "
$ws.Range("E42").Value = "  -0.31%  "
$ws.Range("D43").Value = "'8.33"
$ws.Range("E43").Value = "  -1.88%  "
$ws.Range("D44").Value = "'2.68"
$ws.Range("E44").Value = "  +3.31%  "
$ws.Range("D45").Value = "2.793.56"
$ws.Range("E45").Value = "  -0.59%  "
$ws.Range("D46").Value = "'371.57"
$ws.Range("E46").Value = "  -4.48%  "
$ws.Range("D47").Value = "'135.56"
$ws.Range("E47").Value = "  +0.28%  "
$ws.Range("D48").Value = "'0.0344"
$ws.Range("E48").Value = "  -2.44%  "
$ws.Range("D50").Value = "'24.93"
$ws.Range("E50").Value = "  -0.74%  "
$ws.Range("E51").Value = "  +1.34%  "
